$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.476.81"
$ws.Range("E2").Value = "  -4.31%  "
$ws.Range("D3").Value = "2.375.57"
$ws.Range("E3").Value = "  -5.21%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "500.77"
$ws.Range("E5").Value = "  -6.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.64"
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -3.18%  "
$ws.Range("D9").Value = "2.398.21"
$ws.Range("E9").Value = "  -4.50%  "
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.321"
$ws.Range("E12").Value = "  -3.07%  "
$ws.Range("E13").Value = "  -10.36%  "
$ws.Range("D14").Value = "2.793.54"
$ws.Range("E14").Value = "  -5.29%  "
$ws.Range("D15").Value = "56.947.46"
$ws.Range("E15").Value = "  -3.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.59"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -3.67%  "
$ws.Range("D18").Value = "2.416.26"
$ws.Range("E18").Value = "  -3.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.18"
$ws.Range("E19").Value = "  -4.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "312.49"
$ws.Range("E20").Value = "  -2.85%  "
$ws.Range("E21").Value = "  -5.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.20"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.79"
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "2.490.58"
$ws.Range("E26").Value = "  -4.91%  "
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.371"
$ws.Range("E27").Value = "  -9.33%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.150"
$ws.Range("E28").Value = "  -5.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.21"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0713"
$ws.Range("E31").Value = "  -6.00%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.66"
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.12"
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.10"
$ws.Range("E34").Value = "  -6.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.995"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.79"
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.23"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("E39").Value = "  -5.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "35.82"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.43"
$ws.Range("E41").Value = "  -6.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.778"
$ws.Range("E42").Value = "  -6.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "130.01"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.36"
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.88"
$ws.Range("E45").Value = "  -3.11%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.570"
$ws.Range("E46").Value = "  -4.11%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "253.69"
$ws.Range("E47").Value = "  -8.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0898"
$ws.Range("E48").Value = "  -3.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0486"
$ws.Range("E49").Value = "  -4.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.77"
$ws.Range("E50").Value = "  -4.92%  "
$ws.Range("E51").Value = "  -5.17%  "
